$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Menu > View : add a new "indent width" submenu entry right after
#    the existing "indent type ... submenu" entry.
# ---------------------------------------------------------------------
$anchor = $d.Content.Find.Execute("indent type", $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "", 0)

$indentTypePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "indent type*submenu*") {
        $indentTypePara = $para
        break
    }
}

$srcRange = $indentTypePara.Range
$srcRange.InsertParagraphAfter()

$newIndex = $indentTypePara.Index + 1
$newPara = $d.Paragraphs.Item($newIndex)
$newRange = $newPara.Range
$newRange.Collapse(1)
$newRange.InsertAfter("indent width")
$newRange.Collapse(0)
$newRange.InsertAfter("`t`t`t`t")
$newRange.Collapse(0)
$newRange.InsertAfter(" submenu")

# ---------------------------------------------------------------------
# 2) Menu > Help : remove the first (now-redundant) blank paragraph
#    that immediately follows the "about" entry.
# ---------------------------------------------------------------------
$aboutPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "^\s*about\s*$") {
        $aboutPara = $para
        break
    }
}

$blankIndex = $aboutPara.Index + 1
$blankPara = $d.Paragraphs.Item($blankIndex)
$blankPara.Range.Delete()
